# Add four new "benefit for the vendors & user" statement bullets to the
# end of the existing numbered/bulleted list in the document, matching
# the formatting (style, numbering, spacing, font) of the existing items.

$d = $word.ActiveDocument

$newStatements = @(
    "How can the students maximize their break or lunch time in terms of purchasing food?",
    "How can the vendors accommodate more customers at a time?",
    "What can be a more convenient way to do transactions in terms of payment?",
    "What can help with the inventory process which the vendors do at the end of the day?"
)

foreach ($statement in $newStatements) {
    # Grab the current last paragraph in the document (the previous list
    # item) and collapse its range to the end so we insert right after it.
    $lastPara = $d.Paragraphs.Last
    $r = $lastPara.Range
    $r.Collapse(0)

    # Inserting a new paragraph after it duplicates the paragraph (and
    # run) formatting of that last paragraph - i.e. the ListParagraph
    # style, the list numbering (ilvl 0 / numId 1), the 480 auto line
    # spacing, and the Times New Roman 12pt (sz 24) run formatting.
    $r.InsertParagraphAfter()

    # The freshly inserted paragraph is now the new last paragraph;
    # set its text to the new statement.
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $statement
}
